# Apply cryptos list price/volume updates for Thu Nov  2 17:59:22 UTC 2023
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Cell, $Text)
    # Plain decimal-looking strings (e.g. "0.316") get auto-coerced to
    # numbers by Excel; force text storage, then drop the number-format
    # override so the cell keeps its original (unstyled) appearance.
    $Cell.NumberFormat = "@"
    $Cell.Value = $Text
    $Cell.ClearFormats()
}

$ws.Range("D2").Value = "34.905.55"
$ws.Range("E2").Value = "  +1.07%  "
$ws.Range("D3").Value = "1.811.51"
$ws.Range("E3").Value = "  +0.29%  "
$ws.Range("E4").Value = "  +0.45%  "
Set-TextValue $ws.Range("D5") "232.15"
$ws.Range("E5").Value = "  +3.24%  "
$ws.Range("E6").Value = "  +0.81%  "
$ws.Range("E7").Value = "  +0.29%  "
Set-TextValue $ws.Range("D8") "40.54"
$ws.Range("E8").Value = "  -6.02%  "
Set-TextValue $ws.Range("D9") "0.316"
$ws.Range("E9").Value = "  +7.98%  "
Set-TextValue $ws.Range("D10") "0.0685"
$ws.Range("E10").Value = "  +2.53%  "
Set-TextValue $ws.Range("D11") "0.0997"
$ws.Range("E11").Value = "  -0.05%  "
$ws.Range("D12").Value = "2.072.54"
$ws.Range("E12").Value = "  +0.29%  "
$ws.Range("D13").Value = "1.815.04"
$ws.Range("E13").Value = "  +0.63%  "
$ws.Range("E14").Value = "  +1.86%  "
$ws.Range("E15").Value = "  +6.42%  "
Set-TextValue $ws.Range("D16") "0.658"
$ws.Range("E16").Value = "  +4.41%  "
$ws.Range("D17").Value = "34.865.31"
Set-TextValue $ws.Range("D18") "69.17"
$ws.Range("E18").Value = "  +2.74%  "
$ws.Range("D19").Value = "0.0₃0785"
$ws.Range("E19").Value = "  +2.03%  "
Set-TextValue $ws.Range("D20") "236.92"
$ws.Range("E20").Value = "  -1.56%  "
$ws.Range("E21").Value = "  +5.66%  "
Set-TextValue $ws.Range("D22") "4.65"
$ws.Range("E22").Value = "  +6.37%  "
Set-TextValue $ws.Range("D23") "1.00"
$ws.Range("E23").Value = "  +0.36%  "
Set-TextValue $ws.Range("D24") "2.28"
$ws.Range("E24").Value = "  +6.22%  "
Set-TextValue $ws.Range("D25") "172.53"
$ws.Range("E25").Value = "  +1.21%  "
$ws.Range("E26").Value = "  +2.50%  "
Set-TextValue $ws.Range("D27") "17.43"
$ws.Range("E27").Value = "  -0.10%  "
$ws.Range("E28").Value = "  -0.50%  "
$ws.Range("E29").Value = "  +31.80%  "
$ws.Range("E30").Value = "  +0.65%  "
$ws.Range("D31").Value = "3.339.00"
Set-TextValue $ws.Range("D32") "0.0553"
$ws.Range("E32").Value = "  +7.63%  "
$ws.Range("E33").Value = "  +2.15%  "
Set-TextValue $ws.Range("D34") "3.96"
$ws.Range("E34").Value = "  +2.03%  "
$ws.Range("E35").Value = "  -0.37%  "
Set-TextValue $ws.Range("D36") "93.18"
$ws.Range("E36").Value = "  +5.97%  "
Set-TextValue $ws.Range("D37") "1.14"
$ws.Range("E37").Value = "  +7.25%  "
Set-TextValue $ws.Range("D38") "0.679"
$ws.Range("E38").Value = "  +4.49%  "
Set-TextValue $ws.Range("D39") "0.0193"
$ws.Range("E39").Value = "  +2.31%  "
$ws.Range("D40").Value = "1.307.28"
$ws.Range("E40").Value = "  -0.82%  "
$ws.Range("E41").Value = "  +4.84%  "
$ws.Range("B42").Value = "InjectiveProtocol"
$ws.Range("C42").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue $ws.Range("D42") "14.87"
$ws.Range("E42").Value = "  -0.92%  "
$ws.Range("B43").Value = "ARBITRUM"
$ws.Range("C43").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue $ws.Range("D43") "0.988"
$ws.Range("E43").Value = "  +5.25%  "
$ws.Range("E44").Value = "  -0.66%  "
$ws.Range("E45").Value = "  +0.50%  "
Set-TextValue $ws.Range("D46") "2.76"
$ws.Range("E46").Value = "  -1.38%  "
Set-TextValue $ws.Range("D47") "6.25"
$ws.Range("E47").Value = "  +7.59%  "
$ws.Range("E48").Value = "  -1.20%  "
$ws.Range("D49").Value = "1.988.83"
$ws.Range("E49").Value = "  +1.12%  "
$ws.Range("E50").Value = "  +0.23%  "
$ws.Range("E51").Value = "  +5.49%  "
